$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 116, pushing existing rows 116..181 down to 117..182
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row 116 with the new record's data
$ws.Cells.Item(116, 1).Value = 10
$ws.Cells.Item(116, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(116, 3).Value = "La Araucanía"
$ws.Cells.Item(116, 4).Value = 44606
$ws.Cells.Item(116, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(116, 5).Value = 9
$ws.Cells.Item(116, 6).Value = 100112043
$ws.Cells.Item(116, 7).Value = "Pepino dulce"
$ws.Cells.Item(116, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 270
$ws.Cells.Item(116, 11).Value = 18000
$ws.Cells.Item(116, 12).Value = 19000
$ws.Cells.Item(116, 13).Value = 18444
$ws.Cells.Item(116, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(116, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(116, 16).Value = 1025
$ws.Cells.Item(116, 17).Value = 18
$ws.Cells.Item(116, 18).Value = "Hortaliza"
